$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)
for ($i = $tbl.Rows.Count; $i -ge 2; $i--) {
    $tbl.Rows.Item($i).Delete()
}
